$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 92, shifting rows 92:174 down to 93:175
$ws.Rows(92).Insert()

# Populate the new row 92 with a fresh data entry (same shape as the rows
# around it), changing only the date (D) and volume (J) values.
$ws.Range("A92").Value = 10
$ws.Range("B92").Value = "Vega Modelo de Temuco"
$ws.Range("C92").Value = "La Araucanía"
$ws.Range("D92").Value = 44512
$ws.Range("E92").Value = 9
$ws.Range("F92").Value = 100112039
$ws.Range("G92").Value = "Ciboulette"
$ws.Range("H92").Value = "Sin especificar"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 20
$ws.Range("K92").Value = 5000
$ws.Range("L92").Value = 5000
$ws.Range("M92").Value = 5000
$ws.Range("N92").Value = "$/docena de atados"
$ws.Range("O92").Value = "Provincia de Cautín"
$ws.Range("P92").Value = 1667
$ws.Range("Q92").Value = 3
$ws.Range("R92").Value = "Hortaliza"
